$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "65.331.90"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.172.15"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.07%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "596.39"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.99%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "149.64"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.170.96"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("E9").Value = "  +0.85%  "

# Row 10
$ws.Range("E10").Value = "  -2.41%  "

# Row 11
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.506"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("E13").Value = "  -2.73%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "37.75"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.686.31"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.157.40"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "7.25"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.183.13"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("E19").Value = "  -0.23%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "508.40"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "15.79"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +5.51%  "

# Row 22
$ws.Range("E22").Value = "  -1.80%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "15.13"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.78%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.88"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "84.93"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.18"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "

# Row 28
$ws.Range("E28").Value = "  +2.43%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.22"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30
$ws.Range("E30").Value = "  +2.35%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "27.93"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.68%  "

# Row 32
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.53"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.10%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "6.56"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.97%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "55.11"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0907"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.05%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "476.43"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0421"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.97"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -6.16%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "8.89"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +2.77%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.006.01"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.19%  "

# Row 43
$ws.Range("E43").Value = "  -3.48%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.287"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.76%  "

# Row 46
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0614"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.67%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "28.63"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.13%  "

# Row 49
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.28"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.28%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "119.43"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -3.52%  "
